# Update ICS summary excel in ICS folder
# Adds the 2024 balance row to "Balances" and the corresponding deposit
# (year-over-year change) row to "Deposits", then leaves "Deposits" as
# the active/selected sheet (matching the author's final on-screen state).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Balances!A19:H19 - new "Balance - Dec 2024 (AF)" row
# ---------------------------------------------------------------------
$balances = $wb.Worksheets.Item("Balances")

$balances.Range("A19").Value = "Balance - Dec 2024 (AF)"
$balances.Range("B19").Value = 710589
$balances.Range("C19").Value = 1661832
$balances.Range("D19").Value = 954013
$balances.Range("E19").Value = 240975
$balances.Range("F19").Value = 3326434
$balances.Range("G19").Value = 2024
$balances.Range("H19").Value = 3567409

$balances.Range("B19:F19").NumberFormat = "#,##0"
$balances.Range("H19").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# 2. Deposits!A19:F19 - new 2024 year-over-year deposit row
#    (pick up number formatting from the row above, then fill in the
#    formulas so the stored style matches the existing pattern)
# ---------------------------------------------------------------------
$deposits = $wb.Worksheets.Item("Deposits")

$deposits.Range("A18:F18").Copy()
$deposits.Range("A19").PasteSpecial(-4122)

$deposits.Range("A19").Formula = "=Balances!G19"
$deposits.Range("B19").Formula = "=SUM(C19:F19)"
$deposits.Range("C19").Formula = "=Balances!B19-Balances!B18"
$deposits.Range("D19").Formula = "=Balances!C19-Balances!C18"
$deposits.Range("E19").Formula = "=Balances!D19-Balances!D18"
$deposits.Range("F19").Formula = "=Balances!E19-Balances!E18"

# ---------------------------------------------------------------------
# 3. Leave selection on Deposits (the tab that ends up active), mirroring
#    the cursor ending up on the freshly-typed row.
# ---------------------------------------------------------------------
$balances.Range("A20").Select()
$deposits.Activate()
$deposits.Range("A18:F19").Select()
